$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 previously held the 2023-01-12 (44890) / "La Ligua" batch and
# rows 6-9 held the 2023-02-01 (44908) / "Provincia de Limarí" batch.
# The edit swaps the two weekly batches between the row blocks while
# keeping the "Calidad" (L) ordering (Especial/Primera/Segunda/Tercera)
# the same in each block.

# New values for rows 2-5 (now the 44908 / Provincia de Limarí batch)
$ws.Range("D2").Value = 44908
$ws.Range("M2").Value = 110
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1400

$ws.Range("D3").Value = 44908
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1200

$ws.Range("D4").Value = 44908
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 5000
$ws.Range("P4").Value = 5000
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1000

$ws.Range("D5").Value = 44908
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 4000
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 4000
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 800

# New values for rows 6-9 (now the 44890 / La Ligua batch)
$ws.Range("D6").Value = 44890
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("R6").Value = "La Ligua"
$ws.Range("S6").Value = 2600

$ws.Range("D7").Value = 44890
$ws.Range("M7").Value = 170
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 11000
$ws.Range("P7").Value = 11000
$ws.Range("R7").Value = "La Ligua"
$ws.Range("S7").Value = 2200

$ws.Range("D8").Value = 44890
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("R8").Value = "La Ligua"
$ws.Range("S8").Value = 1600

$ws.Range("D9").Value = 44890
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 7000
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 7000
$ws.Range("R9").Value = "La Ligua"
$ws.Range("S9").Value = 1400
